$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Business demography" row (row 10) with the latest ONS
# Business Demography release dates.
$ws.Range("C10").Value = "Dec 2023 - Dec 2024 (20/11/25)"
$ws.Range("D10").Value = "Dec 2024 - Dec 2025 (Nov 26)"

# Reflect the author's last active selection in the sheet view.
$ws.Range("B22").Select() | Out-Null
